# Machine_Service_Lookup.xlsx edit script
# Commit: إضافة عمود جديد 'Event ' إلى Card19
#
# This commit makes two related changes that were captured in the diff:
#  1. On "Card19": a brand new "Event " column (M) is introduced (mirroring the
#     layout already used on the neighbouring "Card20" sheet). The header style
#     is copied from the existing header cells, the sheet dimension grows by one
#     column, and the previously-blank cells in D:L (rows 2-12) are back-filled
#     with the literal placeholder text "nan" wherever no real value was present
#     (exactly like the equivalent cells already look on "Card20").
#  2. On "Card20": the stray placeholder text "nan" that used to fill every
#     empty cell (columns D:N, rows 2-12) is cleaned out, and the trailing
#     space is trimmed from the "Serviced by " header in O1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Card19 - add the new "Event " column
# ---------------------------------------------------------------------------
$ws19 = $wb.Worksheets.Item("Card19")

# New header cell with the same look & feel as the other header cells (bold,
# centered horizontally, top-aligned vertically, thin border on all sides).
$ws19.Range("M1").Value = "Event "
$ws19.Range("M1").Font.Bold = $true
$ws19.Range("M1").HorizontalAlignment = -4108
$ws19.Range("M1").VerticalAlignment = -4160
$ws19.Range("M1").Borders.LineStyle = 1

# Back-fill every previously empty data cell in D:L (rows 2-12) with "nan",
# leaving cells that already contain real data (ticks, dates, tonnage values)
# untouched.
$card19NanCells = @(
    "D2","E2","F2","G2","H2","I2","J2","K2","L2",
    "D3","F3","G3","H3","I3","J3","K3",
    "D4","E4","F4","G4","H4","I4","J4","K4","L4",
    "D5","E5","F5","G5","H5","I5","J5","K5","L5",
    "E6","F6","G6","J6","K6",
    "E7","G7","H7","I7","J7","K7",
    "D8","E8","F8","G8","H8","I8","J8","K8","L8",
    "D9","E9","F9","G9","H9","I9","J9","K9","L9",
    "D10","E10","F10","G10","H10","I10","J10","K10","L10",
    "D11","E11","F11","G11","H11","I11","J11","K11","L11",
    "D12","E12","F12","G12","H12","I12","J12","K12","L12"
)
foreach ($addr in $card19NanCells) {
    $ws19.Range($addr).Value = "nan"
}

# New column M (rows 2-12) stays empty, same as the neighboring O column on
# Card20 that also has no data.
$card19EmptyM = @("M2","M3","M4","M5","M6","M7","M8","M9","M10","M11","M12")
foreach ($addr in $card19EmptyM) {
    $ws19.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------------
# 2) Card20 - strip the placeholder "nan" text and trim the O1 header
# ---------------------------------------------------------------------------
$ws20 = $wb.Worksheets.Item("Card20")

$ws20.Range("O1").Value = "Serviced by"

$card20ClearCells = @(
    "D2","E2","F2","G2","H2","I2","J2","K2","L2","M2","N2",
    "D3","E3","F3","G3","H3","I3","J3","K3","L3","M3","N3",
    "D4","E4","F4","G4","H4","I4","J4","K4","L4","M4","N4",
    "D5","H5","I5","J5","K5","M5","N5",
    "E6","F6","G6","J6","K6","M6","N6",
    "E7","G7","H7","I7","J7","K7","M7","N7",
    "D8","E8","F8","G8","H8","I8","J8","K8","L8","M8","N8",
    "D9","E9","F9","G9","H9","I9","J9","K9","L9","M9","N9",
    "D10","E10","F10","G10","H10","I10","J10","K10","L10","M10","N10",
    "D11","E11","F11","G11","H11","I11","J11","K11","L11","M11","N11",
    "D12","E12","F12","G12","H12","I12","J12","K12","L12","M12","N12"
)
foreach ($addr in $card20ClearCells) {
    $ws20.Range($addr).ClearContents()
}
